# Apply the 6-May-2020 19:03 COVID-19 data refresh to the "Pais" sheet:
#   - bump the "last updated" banner in A1
#   - rewrite Casos totales/Nuevos casos/Casos activos/Recuperados/Casos
#     criticos/Muertes hoy/Muertes (cols B:H) for every country whose figures
#     moved, re-sorted descending by Casos totales (col B)
#   - a few countries changed rank against their neighbours (Congo vs Ruanda;
#     Santa Sede vs Islas Turcas y Caicos vs Montserrat vs Seychelles), so the
#     country name in col A is rewritten on those rows as well

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 19:03"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1244302
$ws.Range("C4").Value = 6669
$ws.Range("D4").Value = 202908
$ws.Range("E4").Value = 968527
$ws.Range("F4").Value = 16173
$ws.Range("G4").Value = 596
$ws.Range("H4").Value = 72867

# Row 15: Canada
$ws.Range("B15").Value = 62465
$ws.Range("C15").Value = 419
$ws.Range("D15").Value = 26993
$ws.Range("E15").Value = 31429
$ws.Range("F15").Value = 502
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 4043

# Row 28: Irlanda
$ws.Range("B28").Value = 22248
$ws.Range("C28").Value = 265
$ws.Range("D28").Value = 13386
$ws.Range("E28").Value = 7487
$ws.Range("F28").Value = 93
$ws.Range("G28").Value = 36
$ws.Range("H28").Value = 1375

# Row 56: Marruecos
$ws.Range("B56").Value = 5408
$ws.Range("C56").Value = 189
$ws.Range("D56").Value = 2017
$ws.Range("E56").Value = 3208
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 183

# Row 59: Moldavia
$ws.Range("B59").Value = 4476
$ws.Range("C59").Value = 113
$ws.Range("D59").Value = 1658
$ws.Range("E59").Value = 2675
$ws.Range("F59").Value = 237
$ws.Range("G59").Value = 7
$ws.Range("H59").Value = 143

# Row 134: Congo
$ws.Range("A134").Value = "Congo"
$ws.Range("B134").Value = 264
$ws.Range("C134").Value = 28
$ws.Range("D134").Value = 30
$ws.Range("E134").Value = 224
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 10

# Row 135: Ruanda
$ws.Range("A135").Value = "Ruanda"
$ws.Range("B135").Value = 261
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 129
$ws.Range("E135").Value = 132
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

# Row 188: Laos
$ws.Range("B188").Value = 19
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 10
$ws.Range("E188").Value = 9
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Row 203: Santa Sede
$ws.Range("A203").Value = "Santa Sede"
$ws.Range("B203").Value = 12
$ws.Range("C203").Value = 1
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204: Islas Turcas y Caicos
$ws.Range("A204").Value = "Islas Turcas y Caicos"
$ws.Range("B204").Value = 12
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 6
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 1

# Row 205: Montserrat
$ws.Range("A205").Value = "Montserrat"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 7
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 1
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

# Row 206: Seychelles
$ws.Range("A206").Value = "Seychelles"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
